$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 3 - XGBClassifier
$ws.Range("C3").Value = 0.5208333333333334
$ws.Range("D3").Value = 0.4976429989524637
$ws.Range("E3").Value = 0.5422373153382579

# Row 4 - LogisticRegression
$ws.Range("C4").Value = 0.7030228758169935
$ws.Range("D4").Value = 0.6374139434025399
$ws.Range("E4").Value = 0.7670029775506012
